$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the 4-row "day block" structure/styles from the last existing
# block (rows 29-32: header row + "gisteren" row + "vandaag" row + "hulp" row)
# into the new block (rows 33-36).
$src = $ws.Range("A29:D32")
$dst = $ws.Range("A33:D36")
$src.Copy($dst)

# Row 33: header row for the new day (serial date 45376 = 2024-03-25, the
# next day in the log after row 29's 45373 = 2024-03-22).
$ws.Range("A33").Value() = 45376

# Fill in the new day's answers in the same order the shared-string table
# was originally built (Bruno's column top-to-bottom, then Abbi's column
# bottom-to-top, then Chris's column top-to-bottom).
$ws.Range("B34").Value() = "Vorige week vrijdag ben ik op zoek gegaan naar afbeeldingen en begonnen met het schrijven van teksten."
$ws.Range("B35").Value() = "Vandaag wil ik de huisstyling + content op iedere pagina nalopen/afmaken. Daarnaast ga ik met Chris de login stylen."
$ws.Range("D35").Value() = "Vandaag ga ik verder met het maken van de agenda voor de overzicht pagina's."
$ws.Range("D34").Value() = "Vorige week vrijdag heb ik de tarieven pagina gestyled."
$ws.Range("C34").Value() = "Vorige week vrijdag heb ik alle afbeeldingen geoptimaliseerd, de teksten toegevoegd en begonnen met het stylen van de login pagina."
$ws.Range("C35").Value() = "Vandaag ga ik met Bruno op alle pagina's de huisstyling en content nalopen en de login stylen. En we gaan een blok maken met de medewerkers op de over ons pagina."

# Row 36: "Heb je ergens hulp bij nodig?" answers (unchanged boilerplate text,
# carried over by the copy above - re-assert explicitly to be safe).
$ws.Range("B36").Value() = "Momenteel niet. "
$ws.Range("C36").Value() = "Momenteel niet."
$ws.Range("D36").Value() = "Momenteel niet"

# Row heights to match the "gisteren" (45) / "vandaag" (60) rows elsewhere.
$ws.Rows.Item(34).RowHeight = 45
$ws.Rows.Item(35).RowHeight = 60

# Scroll the view down to the new block, and leave the same active cell
# selection the author ended up with.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("C35").Select()
